$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C64").Value = '[name="Texas"]   ''Suspect.'' Hm. Maybe she chased them out.' + "`n"
$ws.Range("C86").Value = '[name="Emperor"]   And you gotta remember the first rule of Penguin Logistics: ''Don''t sweat the small stuff.''' + "`n"
$ws.Range("C87").Value = '[name="Exusiai"]   Wasn''t the first rule ''live for the party'' yesterday?' + "`n"
$ws.Range("C88").Value = '[name="Croissant"]   The one I heard was ''carpe diem.''' + "`n"
$ws.Range("C95").Value = '[name="Exusiai"]   It''s a box of candies. See, it says ''Victorian Gumdrops'' right here...' + "`n"

$ws.Rows.Item(64).AutoFit()
$ws.Rows.Item(86).AutoFit()
$ws.Rows.Item(87).AutoFit()
$ws.Rows.Item(88).AutoFit()
$ws.Rows.Item(95).AutoFit()
